$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing data but keep formatting (column A / header bold+border style) ---
$ws.Cells.ClearContents()

# --- Insert the 4 new columns introduced by this edit: el_Cl, el_mf_Cl, fg_Cl, fg_mf_Cl ---
$ws.Columns("H").Insert()
$ws.Columns("L").Insert()
$ws.Columns("Q").Insert()
$ws.Columns("V").Insert()

# --- Append the 2 new rows (dichlorobenzene goes at the top later; palmitic acid is appended) ---
$ws.Rows("8:9").Insert()

# --- Column A (and header row) use a bold+bordered+centered style; re-stamp it onto the freshly
#     inserted rows 8:9 so they match the rest of column A exactly ---
$ws.Range("A2").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Write the full header + data block in one shot ---
$data = New-Object "object[,]" 9,24
$data[0,0] = "comp_name"
$data[0,1] = "iupac_name"
$data[0,2] = "molecular_formula"
$data[0,3] = "canonical_smiles"
$data[0,4] = "molecular_weight"
$data[0,5] = "xlogp"
$data[0,6] = "el_C"
$data[0,7] = "el_Cl"
$data[0,8] = "el_H"
$data[0,9] = "el_O"
$data[0,10] = "el_mf_C"
$data[0,11] = "el_mf_Cl"
$data[0,12] = "el_mf_H"
$data[0,13] = "el_mf_O"
$data[0,14] = "fg_C-aliph"
$data[0,15] = "fg_C-arom"
$data[0,16] = "fg_Cl"
$data[0,17] = "fg_alcohol"
$data[0,18] = "fg_carboxyl"
$data[0,19] = "fg_mf_C-aliph"
$data[0,20] = "fg_mf_C-arom"
$data[0,21] = "fg_mf_Cl"
$data[0,22] = "fg_mf_alcohol"
$data[0,23] = "fg_mf_carboxyl"
$data[1,0] = "dichlorobenzene"
$data[1,1] = "1,4-dichlorobenzene"
$data[1,2] = "C6H4Cl2"
$data[1,3] = "C1=CC(=CC=C1Cl)Cl"
$data[1,4] = 147
$data[1,5] = 3.4
$data[1,6] = 6
$data[1,7] = 2
$data[1,8] = 4
$data[1,9] = 0
$data[1,10] = 0.4902448979591837
$data[1,11] = 0.4823129251700681
$data[1,12] = 0.02742857142857143
$data[1,13] = 0
$data[1,14] = 0
$data[1,15] = 6
$data[1,16] = 2
$data[1,17] = 0
$data[1,18] = 0
$data[1,19] = 0
$data[1,20] = 0.517673469387755
$data[1,21] = 0.4823129251700681
$data[1,22] = 0
$data[1,23] = 0
$data[2,0] = "oleic acid"
$data[2,1] = "(z)-octadec-9-enoic acid"
$data[2,2] = "C18H34O2"
$data[2,3] = "CCCCCCCCC=CCCCCCCCC(=O)O"
$data[2,4] = 282.5
$data[2,5] = 6.5
$data[2,6] = 18
$data[2,7] = 0
$data[2,8] = 34
$data[2,9] = 2
$data[2,10] = 0.7653026548672566
$data[2,11] = 0
$data[2,12] = 0.121316814159292
$data[2,13] = 0.1132672566371681
$data[2,14] = 17
$data[2,15] = 0
$data[2,16] = 0
$data[2,17] = 0
$data[2,18] = 1
$data[2,19] = 0.8405345132743363
$data[2,20] = 0
$data[2,21] = 0
$data[2,22] = 0
$data[2,23] = 0.1593522123893805
$data[3,0] = "notvalidcomp"
$data[3,1] = "unidentified"
$data[3,2] = $null
$data[3,3] = $null
$data[3,4] = $null
$data[3,5] = $null
$data[3,6] = $null
$data[3,7] = $null
$data[3,8] = $null
$data[3,9] = $null
$data[3,10] = $null
$data[3,11] = $null
$data[3,12] = $null
$data[3,13] = $null
$data[3,14] = $null
$data[3,15] = $null
$data[3,16] = $null
$data[3,17] = $null
$data[3,18] = $null
$data[3,19] = $null
$data[3,20] = $null
$data[3,21] = $null
$data[3,22] = $null
$data[3,23] = $null
$data[4,0] = "dodecane"
$data[4,1] = "dodecane"
$data[4,2] = "C12H26"
$data[4,3] = "CCCCCCCCCCCC"
$data[4,4] = 170.33
$data[4,5] = 6.1
$data[4,6] = 12
$data[4,7] = 0
$data[4,8] = 26
$data[4,9] = 0
$data[4,10] = 0.846192684788352
$data[4,11] = 0
$data[4,12] = 0.1538660247754359
$data[4,13] = 0
$data[4,14] = 12
$data[4,15] = 0
$data[4,16] = 0
$data[4,17] = 0
$data[4,18] = 0
$data[4,19] = 1.000058709563788
$data[4,20] = 0
$data[4,21] = 0
$data[4,22] = 0
$data[4,23] = 0
$data[5,0] = "naphthalene"
$data[5,1] = "naphthalene"
$data[5,2] = "C10H8"
$data[5,3] = "C1=CC=C2C=CC=CC2=C1"
$data[5,4] = 128.17
$data[5,5] = 3.3
$data[5,6] = 10
$data[5,7] = 0
$data[5,8] = 8
$data[5,9] = 0
$data[5,10] = 0.9371147694468284
$data[5,11] = 0
$data[5,12] = 0.06291643910431459
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 10
$data[5,16] = 0
$data[5,17] = 0
$data[5,18] = 0
$data[5,19] = 0
$data[5,20] = 1.000031208551143
$data[5,21] = 0
$data[5,22] = 0
$data[5,23] = 0
$data[6,0] = "capric acid"
$data[6,1] = "decanoic acid"
$data[6,2] = "C10H20O2"
$data[6,3] = "CCCCCCCCCC(=O)O"
$data[6,4] = 172.26
$data[6,5] = 4.1
$data[6,6] = 10
$data[6,7] = 0
$data[6,8] = 20
$data[6,9] = 2
$data[6,10] = 0.6972599558806455
$data[6,11] = 0
$data[6,12] = 0.1170323928944619
$data[6,13] = 0.1857540926506444
$data[6,14] = 9
$data[6,15] = 0
$data[6,16] = 0
$data[6,17] = 0
$data[6,18] = 1
$data[6,19] = 0.7387147335423198
$data[6,20] = 0
$data[6,21] = 0
$data[6,22] = 0
$data[6,23] = 0.2613317078834321
$data[7,0] = "phenol"
$data[7,1] = "phenol"
$data[7,2] = "C6H6O"
$data[7,3] = "C1=CC=C(C=C1)O"
$data[7,4] = 94.11
$data[7,5] = 1.5
$data[7,6] = 6
$data[7,7] = 0
$data[7,8] = 6
$data[7,9] = 1
$data[7,10] = 0.765763468281798
$data[7,11] = 0
$data[7,12] = 0.06426522154925088
$data[7,13] = 0.1700031877590054
$data[7,14] = 0
$data[7,15] = 6
$data[7,16] = 0
$data[7,17] = 1
$data[7,18] = 0
$data[7,19] = 0
$data[7,20] = 0.8193178195728402
$data[7,21] = 0
$data[7,22] = 0.1807140580172139
$data[7,23] = 0
$data[8,0] = "palmitic acid"
$data[8,1] = "hexadecanoic acid"
$data[8,2] = "C16H32O2"
$data[8,3] = "CCCCCCCCCCCCCCCC(=O)O"
$data[8,4] = 256.42
$data[8,5] = 6.4
$data[8,6] = 16
$data[8,7] = 0
$data[8,8] = 32
$data[8,9] = 2
$data[8,10] = 0.7494579205990172
$data[8,11] = 0
$data[8,12] = 0.125793619842446
$data[8,13] = 0.1247874580765931
$data[8,14] = 15
$data[8,15] = 0
$data[8,16] = 0
$data[8,17] = 0
$data[8,18] = 1
$data[8,19] = 0.8244793697839481
$data[8,20] = 0
$data[8,21] = 0
$data[8,22] = 0
$data[8,23] = 0.1755596287341081

$ws.Range("A1:X9").Value = $data
